$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order matters here: new shared-string entries are appended in the order
# they are first assigned, so we replicate the author's original entry
# order (Notes, then the BoardReader note, then the covid19 keyword,
# then the Twitter keyword/notes, then the sheet-name note).

# For Google Alerts block
$ws.Range("C1").Value = "Notes"
$ws.Range("C1").Font.Bold = $true

# For BoardReader block
$ws.Range("C6").Value = "follows keyword in B2"

# For Google Alerts block (continued)
$ws.Range("B2").Value = "covid19"
$ws.Range("B3").Value = 13

# For BoardReader block (continued)
$ws.Range("B7").Value = 13

# For Twitter block
$ws.Range("B10").Value = "#markinourhearts"
$ws.Range("B11").Value = 500

# For Topic Modelling block
$ws.Range("B14").Value = "#markinourhearts"
$ws.Range("B15").Value = "03-08-20 1432"

# Remove the old "Number of Questions" / "Question 1-4" rows (17-21)
$ws.Range("A17:B21").ClearContents()

# Update selection to match final workbook state
$ws.Range("C16").Select()
